$wb = $excel.ActiveWorkbook
$wsAnaliza = $wb.Worksheets.Item(1)
$wsZarzadzanie = $wb.Worksheets.Item(2)

# On "Analiza Finansowa", the "Rodzaj" (type) column (B) for several subjects
# was re-seeded: rows 29, 39, 42, 62 and 72 move from "Finanse behawioralne"
# to "Zaawansowane strategie inwestycyjne", and the old highlighted fill is
# cleared so the cells fall back to the sheet's default (general) style.
$rowsToReseed = @(29, 39, 42, 62, 72)
foreach ($r in $rowsToReseed) {
    $cell = $wsAnaliza.Range("B$r")
    $cell.Value = "Zaawansowane strategie inwestycyjne"
    $cell.ClearFormats()
}

# Update the current selection on "Zarzadzanie Finansowe" (the active sheet)
# to reflect where editing left off.
[void]$wsZarzadzanie.Range("E44").Select()
